# Standardize Whole Foods item numbers: replace the placeholder sample
# rows with the full Raw Item -> Mapped Item lookup table (spaces/dashes
# normalized to the dash-delimited "DD-DDD-D" format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Raw Item (column A), Mapped Item (column B)
$data = @(
    @(2,  '12-046-2',    '12-046-2'),
    @(3,  '12 046 3',    '12-046-3'),
    @(4,  '130251',      '13-025-1'),
    @(5,  '13 025 126',  '13-025-126'),
    @(6,  '1302514',     '13-025-14'),
    @(7,  '1302515',     '13-025-15'),
    @(8,  '1302517',     '13-025-17'),
    @(9,  '130252',      '13-025-2'),
    @(10, '1302520',     '13-025-20'),
    @(11, '13 025 24',   '13-025-24'),
    @(12, '13-025-25',   '13-025-25'),
    @(13, '13-025-26',   '13-025-26'),
    @(14, '130253',      '13-025-3'),
    @(15, '130254',      '13-025-4'),
    @(16, '130255',      '13-025-5'),
    @(17, '130256',      '13-025-6'),
    @(18, '130257',      '13-025-7'),
    @(19, '130258',      '13-025-8'),
    @(20, '130259',      '13-025-9'),
    @(21, '13-027-1',    '13-027-1'),
    @(22, '13-027-2',    '13-027-2'),
    @(23, '13-027-3',    '13-027-3'),
    @(24, '13-027-4',    '13-027-4'),
    @(25, '170411',      '17-041-1'),
    @(26, '170417',      '17-041-7'),
    @(27, '170512',      '17-051-2'),
    @(28, '6016',        '6-016'),
    @(29, '90311',       '9-031-1'),
    @(30, '90881',       '9-088-1')
)

foreach ($entry in $data) {
    $r = $entry[0]
    $rawItem = $entry[1]
    $mappedItem = $entry[2]

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    # Raw item numbers that are all digits (e.g. "130251") would otherwise
    # be auto-coerced to a numeric cell by Excel; format as Text first so
    # they round-trip as strings, same as the already-dashed entries.
    if ($rawItem -match '^[0-9]+$') {
        $cellA.NumberFormat = "@"
    }
    $cellA.Value = $rawItem
    $cellB.Value = $mappedItem
}
